# Update Hgf-Sdc1 NATMI TPM-derived statistics (ligand/receptor expression,
# specificity and edge weight columns G:T) to values recomputed with the
# new TPM input data, per "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04510733333333333
$ws.Range("H2").Value = 0.135322
$ws.Range("I2").Value = 0.001347508866645585
$ws.Range("J2").Value = 0.001347508866645585
$ws.Range("M2").Value = 2.565830333333333
$ws.Range("N2").Value = 7.697490999999999
$ws.Range("O2").Value = 0.0934185609347503
$ws.Range("P2").Value = 0.0934185609347503
$ws.Range("Q2").Value = 0.1157377641224444
$ws.Range("R2").Value = 1.041639877102
$ws.Range("S2").Value = 0.0001258823391688469
$ws.Range("T2").Value = 0.0001258823391688469
$ws.Range("G3").Value = 0.04510733333333333
$ws.Range("H3").Value = 0.135322
$ws.Range("I3").Value = 0.001347508866645585
$ws.Range("J3").Value = 0.001347508866645585
$ws.Range("O3").Value = 0.3847798091300315
$ws.Range("P3").Value = 0.3847798091300315
$ws.Range("Q3").Value = 0.4767099208397778
$ws.Range("R3").Value = 4.290389287558
$ws.Range("S3").Value = 0.000518494204508913
$ws.Range("T3").Value = 0.000518494204508913
$ws.Range("G4").Value = 0.04510733333333333
$ws.Range("H4").Value = 0.135322
$ws.Range("I4").Value = 0.001347508866645585
$ws.Range("J4").Value = 0.001347508866645585
$ws.Range("M4").Value = 13.68376133333333
$ws.Range("N4").Value = 41.051284
$ws.Range("O4").Value = 0.4982080363333638
$ws.Range("P4").Value = 0.4982080363333638
$ws.Range("Q4").Value = 0.6172379837164444
$ws.Range("R4").Value = 5.555141853448
$ws.Range("S4").Value = 0.0006713397463932933
$ws.Range("T4").Value = 0.0006713397463932933
$ws.Range("G5").Value = 0.04510733333333333
$ws.Range("H5").Value = 0.135322
$ws.Range("I5").Value = 0.001347508866645585
$ws.Range("J5").Value = 0.001347508866645585
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6480206666666667
$ws.Range("N5").Value = 1.944062
$ws.Range("O5").Value = 0.0235935936018545
$ws.Range("P5").Value = 0.0235935936018545
$ws.Range("Q5").Value = 0.02923048421822222
$ws.Range("R5").Value = 0.263074357964
$ws.Range("S5").Value = [double]"3.179257657453147E-05"
$ws.Range("T5").Value = [double]"3.179257657453147E-05"
$ws.Range("I6").Value = 0.3371496619592149
$ws.Range("J6").Value = 0.3371496619592149
$ws.Range("M6").Value = 2.565830333333333
$ws.Range("N6").Value = 7.697490999999999
$ws.Range("O6").Value = 0.0934185609347503
$ws.Range("P6").Value = 0.0934185609347503
$ws.Range("Q6").Value = 28.95784140325111
$ws.Range("R6").Value = 260.62057262926
$ws.Range("S6").Value = 0.03149603623986739
$ws.Range("T6").Value = 0.03149603623986739
$ws.Range("I7").Value = 0.3371496619592149
$ws.Range("J7").Value = 0.3371496619592149
$ws.Range("O7").Value = 0.3847798091300315
$ws.Range("P7").Value = 0.3847798091300315
$ws.Range("S7").Value = 0.1297283825769213
$ws.Range("T7").Value = 0.1297283825769213
$ws.Range("I8").Value = 0.3371496619592149
$ws.Range("J8").Value = 0.3371496619592149
$ws.Range("M8").Value = 13.68376133333333
$ws.Range("N8").Value = 41.051284
$ws.Range("O8").Value = 0.4982080363333638
$ws.Range("P8").Value = 0.4982080363333638
$ws.Range("Q8").Value = 154.4342918324711
$ws.Range("R8").Value = 1389.90862649224
$ws.Range("S8").Value = 0.1679706710351579
$ws.Range("T8").Value = 0.1679706710351579
$ws.Range("I9").Value = 0.3371496619592149
$ws.Range("J9").Value = 0.3371496619592149
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6480206666666667
$ws.Range("N9").Value = 1.944062
$ws.Range("O9").Value = 0.0235935936018545
$ws.Range("P9").Value = 0.0235935936018545
$ws.Range("Q9").Value = 7.313531003035556
$ws.Range("R9").Value = 65.82177902732001
$ws.Range("S9").Value = 0.00795457210726834
$ws.Range("T9").Value = 0.00795457210726834
$ws.Range("G10").Value = 0.8868746666666668
$ws.Range("H10").Value = 2.660624
$ws.Range("I10").Value = 0.02649395095261704
$ws.Range("J10").Value = 0.02649395095261704
$ws.Range("M10").Value = 2.565830333333333
$ws.Range("N10").Value = 7.697490999999999
$ws.Range("O10").Value = 0.0934185609347503
$ws.Range("P10").Value = 0.0934185609347503
$ws.Range("Q10").Value = 2.275569921598222
$ws.Range("R10").Value = 20.480129294384
$ws.Range("S10").Value = 0.002475026771469341
$ws.Range("T10").Value = 0.002475026771469341
$ws.Range("G11").Value = 0.8868746666666668
$ws.Range("H11").Value = 2.660624
$ws.Range("I11").Value = 0.02649395095261704
$ws.Range("J11").Value = 0.02649395095261704
$ws.Range("O11").Value = 0.3847798091300315
$ws.Range("P11").Value = 0.3847798091300315
$ws.Range("Q11").Value = 9.37279863159289
$ws.Range("R11").Value = 84.355187684336
$ws.Range("S11").Value = 0.0101943373906484
$ws.Range("T11").Value = 0.0101943373906484
$ws.Range("G12").Value = 0.8868746666666668
$ws.Range("H12").Value = 2.660624
$ws.Range("I12").Value = 0.02649395095261704
$ws.Range("J12").Value = 0.02649395095261704
$ws.Range("M12").Value = 13.68376133333333
$ws.Range("N12").Value = 41.051284
$ws.Range("O12").Value = 0.4982080363333638
$ws.Range("P12").Value = 0.4982080363333638
$ws.Range("Q12").Value = 12.13578127124623
$ws.Range("R12").Value = 109.222031441216
$ws.Range("S12").Value = 0.01319949927881579
$ws.Range("T12").Value = 0.01319949927881579
$ws.Range("G13").Value = 0.8868746666666668
$ws.Range("H13").Value = 2.660624
$ws.Range("I13").Value = 0.02649395095261704
$ws.Range("J13").Value = 0.02649395095261704
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6480206666666667
$ws.Range("N13").Value = 1.944062
$ws.Range("O13").Value = 0.0235935936018545
$ws.Range("P13").Value = 0.0235935936018545
$ws.Range("Q13").Value = 0.5747131127431112
$ws.Range("R13").Value = 5.172418014688001
$ws.Range("S13").Value = 0.0006250875116835122
$ws.Range("T13").Value = 0.0006250875116835122
$ws.Range("G14").Value = 21.25667433333333
$ws.Range("H14").Value = 63.77002299999999
$ws.Range("I14").Value = 0.6350088782215225
$ws.Range("J14").Value = 0.6350088782215224
$ws.Range("M14").Value = 2.565830333333333
$ws.Range("N14").Value = 7.697490999999999
$ws.Range("O14").Value = 0.0934185609347503
$ws.Range("P14").Value = 0.0934185609347503
$ws.Range("Q14").Value = 54.54101979025477
$ws.Range("R14").Value = 490.8691781122929
$ws.Range("S14").Value = 0.05932161558424474
$ws.Range("T14").Value = 0.05932161558424472
$ws.Range("G15").Value = 21.25667433333333
$ws.Range("H15").Value = 63.77002299999999
$ws.Range("I15").Value = 0.6350088782215225
$ws.Range("J15").Value = 0.6350088782215224
$ws.Range("O15").Value = 0.3847798091300315
$ws.Range("P15").Value = 0.3847798091300315
$ws.Range("Q15").Value = 224.6478962495441
$ws.Range("R15").Value = 2021.831066245897
$ws.Range("S15").Value = 0.2443385949579528
$ws.Range("T15").Value = 0.2443385949579528
$ws.Range("G16").Value = 21.25667433333333
$ws.Range("H16").Value = 63.77002299999999
$ws.Range("I16").Value = 0.6350088782215225
$ws.Range("J16").Value = 0.6350088782215224
$ws.Range("M16").Value = 13.68376133333333
$ws.Range("N16").Value = 41.051284
$ws.Range("O16").Value = 0.4982080363333638
$ws.Range("P16").Value = 0.4982080363333638
$ws.Range("Q16").Value = 290.8712583177258
$ws.Range("R16").Value = 2617.841324859532
$ws.Range("S16").Value = 0.3163665262729969
$ws.Range("T16").Value = 0.3163665262729968
$ws.Range("G17").Value = 21.25667433333333
$ws.Range("H17").Value = 63.77002299999999
$ws.Range("I17").Value = 0.6350088782215225
$ws.Range("J17").Value = 0.6350088782215224
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6480206666666667
$ws.Range("N17").Value = 1.944062
$ws.Range("O17").Value = 0.0235935936018545
$ws.Range("P17").Value = 0.0235935936018545
$ws.Range("Q17").Value = 13.77476427260289
$ws.Range("R17").Value = 123.972878453426
$ws.Range("S17").Value = 0.01498214140632811
$ws.Range("T17").Value = 0.01498214140632811
